# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" (E16:E88) labels are re-pointed to the most recent 73
# periods (2207 down to 1607, newest first) replacing the previous ordering
# (1607 up to 2207, oldest first), and the "Valor Mora" (G16:G88) amounts are
# all raised to 6000000. The first data row's "Salario Basico" (F16) drops to
# 50000 while the former last-row value (F88) moves up to 60000, so column F
# ends up uniform at 60000 for every row except the first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newest-to-oldest replacement list for the "Periodo Mora" column (E16:E88).
$periods = @(
    "2207", "2206", "2205", "2204", "2203", "2202", "2201",
    "2112", "2111", "2110", "2109", "2108", "2107", "2106", "2105", "2104", "2103", "2102", "2101",
    "2012", "2011", "2010", "2009", "2008", "2007", "2006", "2005", "2004", "2003", "2002", "2001",
    "1912", "1911", "1910", "1909", "1908", "1907", "1906", "1905", "1904", "1903", "1902", "1901",
    "1812", "1811", "1810", "1809", "1808", "1807", "1806", "1805", "1804", "1803", "1802", "1801",
    "1712", "1711", "1710", "1709", "1708", "1707", "1706", "1705", "1704", "1703", "1702", "1701",
    "1612", "1611", "1610", "1609", "1608", "1607"
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 5).Value = $periods[$i]
}

# "Valor Mora" (column G) becomes 6000000 for every data row.
$ws.Range("G16:G88").Value = 6000000

# "Salario Basico" (column F): row 16 drops to 50000, row 88 rises to 60000
# (the other rows already hold 60000 and are left untouched).
$ws.Range("F16").Value = 50000
$ws.Range("F88").Value = 60000
